$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "server3" column header (I2) ---
# Start from the existing "server1" header (G2) which already has the
# font/fill/border/alignment combination we want (bold font + themed
# fill + thin left/right/top/bottom border, centered), then strip the
# top/bottom border edges so it matches the new header style (left/right
# border only).
$ws.Range("I2").Value = "server3"
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Borders.Item(8).LineStyle = -4142
$ws.Range("I2").Borders.Item(9).LineStyle = -4142

# --- New "server3" IP addresses (I3:I16) ---
# Copy the format from the matching H-column data cell (same row) so the
# new column's formatting (font/border) matches the existing G/H columns.
$values = @("13.127.183.236","3.6.94.118","13.201.63.178","43.205.208.81","3.110.214.199","13.201.90.70","65.1.135.192","13.126.242.36","13.233.17.29","15.206.68.190","65.1.1.83","15.206.157.147","13.232.142.149","65.0.17.37")

for ($r = 3; $r -le 16; $r++) {
    $idx = $r - 3
    $ws.Range("H$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
    $ws.Range("I$r").Value = $values[$idx]
}

# Match the saved selection state from the authored workbook.
$ws.Range("I17").Select()
